$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> FAPs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Gnai2"
$ws.Cells.Item(2,3).Value = "Cnr1"
$ws.Cells.Item(2,4).Value = "FAPs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 63.91118233333333
$ws.Cells.Item(2,8).Value = 191.733547
$ws.Cells.Item(2,9).Value = 0.4067926910433548
$ws.Cells.Item(2,10).Value = 0.4067926910433549
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 1.220413333333333
$ws.Cells.Item(2,14).Value = 3.66124
$ws.Cells.Item(2,15).Value = 0.9446596300458959
$ws.Cells.Item(2,16).Value = 0.9446596300458959
$ws.Cells.Item(2,17).Value = 77.99805906869777
$ws.Cells.Item(2,18).Value = 701.9825316182799
$ws.Cells.Item(2,19).Value = 0.38428063302639
$ws.Cells.Item(2,20).Value = 0.38428063302639

# Row 3: ECs -> MuSCs
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Gnai2"
$ws.Cells.Item(3,3).Value = "Cnr1"
$ws.Cells.Item(3,4).Value = "MuSCs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 63.91118233333333
$ws.Cells.Item(3,8).Value = 191.733547
$ws.Cells.Item(3,9).Value = 0.4067926910433548
$ws.Cells.Item(3,10).Value = 0.4067926910433549
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 0.6666666666666666
$ws.Cells.Item(3,13).Value = 0.07149466666666666
$ws.Cells.Item(3,14).Value = 0.214484
$ws.Cells.Item(3,15).Value = 0.05534036995410405
$ws.Cells.Item(3,16).Value = 0.05534036995410405
$ws.Cells.Item(3,17).Value = 4.569308677194222
$ws.Cells.Item(3,18).Value = 41.123778094748
$ws.Cells.Item(3,19).Value = 0.02251205801696481
$ws.Cells.Item(3,20).Value = 0.02251205801696481

# Row 4: FAPs -> FAPs
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Gnai2"
$ws.Cells.Item(4,3).Value = "Cnr1"
$ws.Cells.Item(4,4).Value = "FAPs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 57.4434
$ws.Cells.Item(4,8).Value = 172.3302
$ws.Cells.Item(4,9).Value = 0.3656254573230189
$ws.Cells.Item(4,10).Value = 0.365625457323019
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 1.220413333333333
$ws.Cells.Item(4,14).Value = 3.66124
$ws.Cells.Item(4,15).Value = 0.9446596300458959
$ws.Cells.Item(4,16).Value = 0.9446596300458959
$ws.Cells.Item(4,17).Value = 70.104691272
$ws.Cells.Item(4,18).Value = 630.9422214480001
$ws.Cells.Item(4,19).Value = 0.3453916092501246
$ws.Cells.Item(4,20).Value = 0.3453916092501246

# Row 5: FAPs -> MuSCs
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Gnai2"
$ws.Cells.Item(5,3).Value = "Cnr1"
$ws.Cells.Item(5,4).Value = "MuSCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 57.4434
$ws.Cells.Item(5,8).Value = 172.3302
$ws.Cells.Item(5,9).Value = 0.3656254573230189
$ws.Cells.Item(5,10).Value = 0.365625457323019
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.07149466666666666
$ws.Cells.Item(5,14).Value = 0.214484
$ws.Cells.Item(5,15).Value = 0.05534036995410405
$ws.Cells.Item(5,16).Value = 0.05534036995410405
$ws.Cells.Item(5,17).Value = 4.106896735199999
$ws.Cells.Item(5,18).Value = 36.9620706168
$ws.Cells.Item(5,19).Value = 0.02023384807289435
$ws.Cells.Item(5,20).Value = 0.02023384807289435

# Row 6: MuSCs -> FAPs
$ws.Cells.Item(6,1).Value = "MuSCs"
$ws.Cells.Item(6,2).Value = "Gnai2"
$ws.Cells.Item(6,3).Value = "Cnr1"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 35.755375
$ws.Cells.Item(6,8).Value = 107.266125
$ws.Cells.Item(6,9).Value = 0.2275818516336261
$ws.Cells.Item(6,10).Value = 0.2275818516336262
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 1.220413333333333
$ws.Cells.Item(6,14).Value = 3.66124
$ws.Cells.Item(6,15).Value = 0.9446596300458959
$ws.Cells.Item(6,16).Value = 0.9446596300458959
$ws.Cells.Item(6,17).Value = 43.63633638833333
$ws.Cells.Item(6,18).Value = 392.727027495
$ws.Cells.Item(6,19).Value = 0.2149873877693812
$ws.Cells.Item(6,20).Value = 0.2149873877693813

# Row 7: MuSCs -> MuSCs
$ws.Cells.Item(7,1).Value = "MuSCs"
$ws.Cells.Item(7,2).Value = "Gnai2"
$ws.Cells.Item(7,3).Value = "Cnr1"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 35.755375
$ws.Cells.Item(7,8).Value = 107.266125
$ws.Cells.Item(7,9).Value = 0.2275818516336261
$ws.Cells.Item(7,10).Value = 0.2275818516336262
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 0.6666666666666666
$ws.Cells.Item(7,13).Value = 0.07149466666666666
$ws.Cells.Item(7,14).Value = 0.214484
$ws.Cells.Item(7,15).Value = 0.05534036995410405
$ws.Cells.Item(7,16).Value = 0.05534036995410405
$ws.Cells.Item(7,17).Value = 2.556318617166667
$ws.Cells.Item(7,18).Value = 23.0068675545
$ws.Cells.Item(7,19).Value = 0.01259446386424489
$ws.Cells.Item(7,20).Value = 0.01259446386424489

# Remove rows 8-10 (extra Sending=MuSCs / Target=ECs combos no longer present)
$ws.Rows.Item(8).EntireRow.Delete() | Out-Null
$ws.Rows.Item(8).EntireRow.Delete() | Out-Null
$ws.Rows.Item(8).EntireRow.Delete() | Out-Null
